$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 45, pushing existing rows 45+ down to 47+
$ws.Rows.Item(45).Resize(2).Insert()

# New row 45
$ws.Cells.Item(45, 1).Value = 6
$ws.Cells.Item(45, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(45, 3).Value = "Metropolitana"
$ws.Cells.Item(45, 4).Value = 44708
$ws.Cells.Item(45, 5).Value = 13
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100107
$ws.Cells.Item(45, 8).Value = "Otros"
$ws.Cells.Item(45, 9).Value = 100107001
$ws.Cells.Item(45, 10).Value = "Caqui"
$ws.Cells.Item(45, 11).Value = "Mankaki"
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 20
$ws.Cells.Item(45, 14).Value = 280000
$ws.Cells.Item(45, 15).Value = 280000
$ws.Cells.Item(45, 16).Value = 280000
$ws.Cells.Item(45, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(45, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(45, 19).Value = 622
$ws.Cells.Item(45, 20).Value = 450

# New row 46
$ws.Cells.Item(46, 1).Value = 6
$ws.Cells.Item(46, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(46, 3).Value = "Metropolitana"
$ws.Cells.Item(46, 4).Value = 44708
$ws.Cells.Item(46, 5).Value = 13
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100107
$ws.Cells.Item(46, 8).Value = "Otros"
$ws.Cells.Item(46, 9).Value = 100107001
$ws.Cells.Item(46, 10).Value = "Caqui"
$ws.Cells.Item(46, 11).Value = "Mankaki"
$ws.Cells.Item(46, 12).Value = "Segunda"
$ws.Cells.Item(46, 13).Value = 15
$ws.Cells.Item(46, 14).Value = 230000
$ws.Cells.Item(46, 15).Value = 230000
$ws.Cells.Item(46, 16).Value = 230000
$ws.Cells.Item(46, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(46, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(46, 19).Value = 511
$ws.Cells.Item(46, 20).Value = 450
